# Add a new "Jenis tipe soal (untuk essay terbatas/bebas)" question-type
# column (L) to the question-import template, mirroring the header style
# used by the existing columns, widening column L, and moving the
# selection/scroll position over to the newly added column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell with the same look (font/fill/alignment/wrap) as the
# other header cells in row 1 (copy format from K1, the previous last col).
$ws.Range("L1").Value = "Jenis tipe soal`n(untuk essay`nterbatas/bebas)"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen the new column to fit the wrapped header text.
$ws.Columns("L").ColumnWidth = 25

# Move the selection to the row below the new header (mirrors the
# original sheet, which was left with the selection on row 4) and bring
# the new column into view.
$ws.Range("J1").Select()
$ws.Range("L4").Select()
